# regen sval data to filter save games
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data values (B:E) for rows 2-5; G is the sum of B..E (the "sum" column).
$data = @{
    2 = @(3.286832544864788, 1.655778082260271, 22.3905356188092, 0.4942365360607697)
    3 = @(0.6606524410359556, 1.655778082260271, 0.7527432677738641, 0.4942365360607697)
    4 = @(3.286832544864788, 1.655778082260271, 0.7527432677738641, 0.4942365360607697)
    5 = @(3.286832544864788, 1.655778082260271, 0.7527432677738641, 0.4942365360607697)
}

$sums = @{
    2 = 27.82738278199502
    3 = 3.56341032713086
    4 = 6.189590430959694
    5 = 6.189590430959694
}

foreach ($row in 2..5) {
    $vals = $data[$row]
    $ws.Cells.Item($row, 2).Value = $vals[0]  # B
    $ws.Cells.Item($row, 3).Value = $vals[1]  # C
    $ws.Cells.Item($row, 4).Value = $vals[2]  # D
    $ws.Cells.Item($row, 5).Value = $vals[3]  # E
    $ws.Cells.Item($row, 7).Value = $sums[$row]  # G
}
